# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-name suffixes to "_FV2210"/"_FV2304"
# - Wrap the data range in an Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header renames: "_old" columns (A..J) become "_FV2210"; "_new" columns (L..U) become "_FV2304"
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $fv2210Headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2210Headers[$i]
}

for ($i = 0; $i -lt $fv2304Headers.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2304Headers[$i]
}

# 2) Turn the data range into a Table named "Table1"
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U94"), $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Applied AHB regeneration edits"
